$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 160, shifting existing rows 160:193 down to 161:194
$ws.Rows("160").Insert()

# Populate the newly inserted row 160 with the new data point
$ws.Range("A160").Value = 11
$ws.Range("B160").Value = "Vega Monumental Concepción"
$ws.Range("C160").Value = "Bíobío"
$ws.Range("D160").Value = 44782
$ws.Range("E160").Value = 8
$ws.Range("F160").Value = "Fruta"
$ws.Range("G160").Value = 100101
$ws.Range("H160").Value = "Berries"
$ws.Range("I160").Value = 100101007
$ws.Range("J160").Value = "Kiwi"
$ws.Range("K160").Value = "Hayward"
$ws.Range("L160").Value = "Primera"
$ws.Range("M160").Value = 270
$ws.Range("N160").Value = 7000
$ws.Range("O160").Value = 7500
$ws.Range("P160").Value = 7278
$ws.Range("Q160").Value = '$/bandeja 18 kilos'
$ws.Range("R160").Value = "Provincia de Curicó"
$ws.Range("S160").Value = 404
$ws.Range("T160").Value = 18
